$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(8, 4, 40617.33, 13, -1, 0.9),
    @(17, 4, 46440.69, 16, -1, 0.9500000000000001),
    @(21, 4, 33532.21, 10, -1, 0.821875),
    @(27, 1, 8027.79, 3, -1, 0.215625),
    @(46, 4, 29952.92, 11, -1, 0.7875),
    @(58, 4, 26792.33, 9, -1, 0.6687500000000001),
    @(69, 4, 7697.42, 5, -1, 0.371875),
    @(70, 1, 11287.01, 3, -1, 0.246875),
    @(80, 4, 11664.86, 4, -1, 0.4375),
    @(85, 1, 13238.11, 3, -1, 0.340625),
    @(97, 4, 26064.56, 10, -1, 0.665625),
    @(126, 1, 2482.48, 2, -1, 0.096875),
    @(156, 4, 32901.47, 6, -1, 0.74375),
    @(157, 1, 12578.63, 2, -1, 0.253125),
    @(158, 3, 28672.08, 5, -1, 0.5906250000000001),
    @(159, 2, 21572.56, 3, -1, 0.409375)
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $r++
}
